# The sheet is a route_no -> zipcode lookup table (columns B and C).
# A new zipcode "J0N" was added under the existing "YUL" route, inserted
# right before the previous row 470 (B470 was "YUL"/C470 was "J7Z"),
# pushing all the rows below it down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at 470 (existing row 470 and everything below
# shifts down to 471+).
$ws.Rows.Item(470).Insert()

# Fill in the new row's route_no (B) and zipcode (C) values.
$ws.Cells.Item(470, 2).Value = "YUL"
$ws.Cells.Item(470, 3).Value = "J0N"

# Reflect the view/selection state left behind by the edit (scrolled so
# row 448 is at the top, with the newly-typed cell B470 selected).
$ws.Application.Goto($ws.Range("A448"))
$excel.ActiveWindow.ScrollRow = 448
$ws.Range("B470").Select() | Out-Null
